$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column (H1): copy the exact formatting of the
# neighboring header cell (G1) so it shares the same cell style, then set
# the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column (rows 2-7)
$saveValues = @(0, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
